{"js": "// The document contains the sequence of runs:\n//   \"<id>\"  (Courier New, color 7f6000, sz 18)\n//   \"p161r_1\"  (default formatting)\n//   \"</id>\"  (Courier New, color 7f6000, sz 18)\n// which together render the text \"<id>p161r_1</id>\". The edit merges\n// these three runs into a single run carrying the \"<id>...</id>\" tag\n// formatting (Courier New / 7f6000 / 9pt), i.e. the plain run in the\n// middle picks up the tag styling instead of being its own run.\nconst body = context.document.body;\n\n// Locate the exact text span (it reads back correctly even though it is\n// currently split across three runs/formats).\nconst results = body.search(\"<id>p161r_1</id>\", { matchCase: true, ignoreSpace: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find '<id>p161r_1</id>' in the document\");\n}\n\nconst target = results.items[0];\n\n// Replacing the range's text with itself collapses the underlying runs\n// into a single run, and that run inherits the formatting of the first\n// run of the original range (the \"<id>\" tag's Courier New / 7f6000 style).\ntarget.insertText(\"<id>p161r_1</id>\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# The document contains the sequence of runs:\n#   \"<id>\"      (Courier New, color 7f6000, sz 18)\n#   \"p161r_1\"   (default formatting)\n#   \"</id>\"     (Courier New, color 7f6000, sz 18)\n# which together render the text \"<id>p161r_1</id>\". The edit merges\n# these three runs into a single run carrying the \"<id>...</id>\" tag\n# formatting (Courier New / 7f6000 / 9pt).\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"<id>p161r_1</id>\"\n$find.Replacement.Text = \"<id>p161r_1</id>\"\n\n# A Find & Replace over the whole matched span collapses it into a single\n# run that inherits the formatting in effect at the start of the match\n# (the \"<id>\" tag's Courier New / 7f6000 style).\n$find.Execute(\n    $find.Text,\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    $find.Replacement.Text,\n    2\n) | Out-Null\n"}
